# Update feature_codebook worksheet to match the revised codebook:
# - remove the IDU, in_simulation, MMT and incarceration variable blocks
# - expand ART / PrEP from single boolean rows into 3-row categorical blocks
#   (not on / adherent / not adherent), with PrEP now listed before ART

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously used range (A1:D22) first so that rows beyond the
# new extent (A1:D20) don't retain stale values.
$ws.Range("A1:D22").ClearContents()

# Header
$ws.Range("A1").Value = "variable_name"
$ws.Range("B1").Value = "type"
$ws.Range("C1").Value = "values"
$ws.Range("D1").Value = "description"

# age
$ws.Range("A3").Value = "age"
$ws.Range("B3").Value = "integer"
$ws.Range("C3").Value = "0-100"
$ws.Range("D3").Value = "age of person rounded to nearest integer"

# race (categorical)
$ws.Range("A4").Value = "race"
$ws.Range("B4").Value = "categorical"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = "white"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "black"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = "latinX"
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = "asian/pacific islander"
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = "other"

# MSM (boolean)
$ws.Range("A9").Value = "MSM"
$ws.Range("B9").Value = "boolean"
$ws.Range("C9").Value = "0/1"
$ws.Range("D9").Value = "men who have sex with men"

# alive (boolean)
$ws.Range("A10").Value = "alive"
$ws.Range("B10").Value = "boolean"
$ws.Range("C10").Value = "0/1"
$ws.Range("D10").Value = "alive in the current time period"

# HIV (categorical)
$ws.Range("A11").Value = "HIV"
$ws.Range("B11").Value = "categorical"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = "no HIV"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = "early stage HIV"
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = "late stage HIV"
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = "AIDS"

# ART (now categorical: not on / adherent / not adherent) - written before the
# PrEP block below so new shared strings land in the same append order as the
# source workbook (ART descriptions precede PrEP descriptions in sharedStrings.xml).
$ws.Range("A18").Value = "ART"
$ws.Range("B18").Value = "categorical"
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = "not on ART (anti-retroviral therapy)"
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = "on ART, adherent"
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = "was on ART, not adherent"

# PrEP (now categorical: not on / adherent / not adherent)
$ws.Range("A15").Value = "PrEP"
$ws.Range("B15").Value = "boolean"
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = "not on PrEP (pre-exposure prophylaxis)"
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = "on PrEP, adherent"
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = "on PrEP, not adherent"

# Match the author's post-edit selection/scroll state.
$ws.Range("D18").Select()
